$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Added additional parameter validation: reorder the "Force of infection" (foi)
# parameter row so it sits after the other computed-rate parameters instead of
# before them. This moves row 5 (foi) down to the bottom of the parameter list
# (row 8), shifting recrate / infdeath / susdeath up by one row each.

# Remove the "foi" row; recrate/infdeath/susdeath (and their formatting) shift up.
$ws.Rows.Item(5).Delete()

# Re-create the "foi" row at the end of the table (now the blank row 8).
$ws.Cells.Item(8, 1).Value2 = "foi"
$ws.Cells.Item(8, 2).Value2 = "Force of infection"
$ws.Cells.Item(8, 3).Value2 = "Probability"
$ws.Cells.Item(8, 6).Value2 = "(1 - (1-ch_prev*transpercontact)**floor(contacts)*(1-ch_prev*transpercontact*(contacts-floor(contacts))))*(1-susdeath)"

# Match the original formatting of that row (left-aligned code/name columns,
# centered blank "Databook Page" cell, and the "Function" column's font).
$ws.Cells.Item(8, 1).HorizontalAlignment = -4131
$ws.Cells.Item(8, 2).HorizontalAlignment = -4131
$ws.Cells.Item(8, 4).HorizontalAlignment = -4108
$ws.Cells.Item(8, 6).Font.Name = "Calibri"

# Leave the sheet with the moved row selected, matching the post-move UI state.
$ws.Range("A5:XFD5").Select() | Out-Null
